# Rename TestObject sheets to use the "V_" naming convention, separating
# them from the functional test sheets, and update the saved view state
# (active tab, selection, scroll position) to match what the author left
# the workbook in after making the edit.

$wb = $excel.ActiveWorkbook

$wsHome            = $wb.Worksheets.Item("Home_page")
$wsAboutBento      = $wb.Worksheets.Item("AboutBentoPage")
$wsAboutResources  = $wb.Worksheets.Item("AboutResourcesPage")

# --- Rename sheets ---------------------------------------------------------
$wsHome.Name           = "V_HomePage"
$wsAboutBento.Name     = "V_AboutBentoPage"
$wsAboutResources.Name = "V_AboutResourcesPage"

# --- V_HomePage view: no longer the selected tab, scrolled down a bit,
#     selection moved to A40 -------------------------------------------------
$wsHome.Activate() | Out-Null
$wsHome.Range("A40").Select() | Out-Null
$excel.ActiveWindow.ScrollRow    = 3
$excel.ActiveWindow.ScrollColumn = 1

# --- V_AboutBentoPage view: becomes the active/selected tab, selection
#     moved to A6 ------------------------------------------------------------
$wsAboutBento.Activate() | Out-Null
$wsAboutBento.Range("A6").Select() | Out-Null

# --- Workbook window position (best effort) ---------------------------------
$w1 = $wb.Windows.Item(1)
$w1.Left = 8840
$w1.Top  = 1260
